# Updated cryptos list on Tue Jun 20 13:28:42 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row. Values prefixed with a leading "'" are numeric-looking strings
# (e.g. "0.9994") that must stay stored as text, matching the source data
# feed which always writes these columns as plain strings rather than
# numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.872.19'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '1.729.95'
$ws.Range('D4').Value = '''0.9994'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''240.12'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = '''0.4832'
$ws.Range('E7').Value = '  -1.27%  '
$ws.Range('D8').Value = '''0.2595'
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('D9').Value = '''0.06174'
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('D10').Value = '1.729.25'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').Value = '''16.02'
$ws.Range('E11').Value = '  +2.55%  '
$ws.Range('D12').Value = '''0.06873'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').Value = '''0.6036'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('D14').Value = '''4.461'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').Value = '''1.000'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '26.655.21'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '''0.9995'
$ws.Range('D19').Value = '''0.000007116'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').Value = '''11.36'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = '1.951.08'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '''4.394'
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('D23').Value = '''8.422'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').Value = '''5.058'
$ws.Range('E24').Value = '  -0.99%  '
$ws.Range('D25').Value = '''139.89'
$ws.Range('E25').Value = '  +1.77%  '
$ws.Range('D26').Value = '''15.20'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').Value = '''1.796'
$ws.Range('D28').Value = '''106.68'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '''1.380'
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('D30').Value = '''3.956'
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').Value = '''0.07916'
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').Value = '''0.04583'
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('D34').Value = '''2.591'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('D35').Value = '''0.9997'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').Value = '''0.6166'
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('D37').Value = '''0.9244'
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('D38').Value = '''2.458'
$ws.Range('E38').Value = '  +2.98%  '
$ws.Range('D39').Value = '''1.988'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').Value = '''0.9990'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').Value = '''5.708'
$ws.Range('E41').Value = '  +5.48%  '
$ws.Range('D42').Value = '''0.01497'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('D43').Value = '''99.91'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '''0.3834'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').Value = '''6.768'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('D47').Value = '''0.05364'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').Value = '''7.868'
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('D49').Value = '''30.07'
$ws.Range('E49').Value = '  -1.51%  '
$ws.Range('D50').Value = '''1.239'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('D51').Value = '''51.27'
$ws.Range('E51').Value = '  -0.33%  '
